# Replace the placeholder "word" column (B2:B193) with the real working
# set of German verb infinitives for this retrieval sequence, as described
# by the commit "Add working set of sequences". All other columns/data are
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'leugnen'
$ws.Range("B3").Value = 'hören'
$ws.Range("B4").Value = 'greifen'
$ws.Range("B5").Value = 'flüchten'
$ws.Range("B6").Value = 'spielen'
$ws.Range("B7").Value = 'scheitern'
$ws.Range("B8").Value = 'sperren'
$ws.Range("B9").Value = 'lohnen'
$ws.Range("B10").Value = 'weichen'
$ws.Range("B11").Value = 'landen'
$ws.Range("B12").Value = 'enden'
$ws.Range("B13").Value = 'schwächen'
$ws.Range("B14").Value = 'schulden'
$ws.Range("B15").Value = 'geben'
$ws.Range("B16").Value = 'hindern'
$ws.Range("B17").Value = 'freuen'
$ws.Range("B18").Value = 'loben'
$ws.Range("B19").Value = 'machen'
$ws.Range("B20").Value = 'ahnen'
$ws.Range("B21").Value = 'münzen'
$ws.Range("B22").Value = 'spüren'
$ws.Range("B23").Value = 'zögern'
$ws.Range("B24").Value = 'bitten'
$ws.Range("B25").Value = 'sprengen'
$ws.Range("B26").Value = 'siegen'
$ws.Range("B27").Value = 'kürzen'
$ws.Range("B28").Value = 'betteln'
$ws.Range("B29").Value = 'runden'
$ws.Range("B30").Value = 'trennen'
$ws.Range("B31").Value = 'stechen'
$ws.Range("B32").Value = 'fließen'
$ws.Range("B33").Value = 'jubeln'
$ws.Range("B34").Value = 'schleppen'
$ws.Range("B35").Value = 'kehren'
$ws.Range("B36").Value = 'melden'
$ws.Range("B37").Value = 'lesen'
$ws.Range("B38").Value = 'klingeln'
$ws.Range("B39").Value = 'schwören'
$ws.Range("B40").Value = 'flehen'
$ws.Range("B41").Value = 'biegen'
$ws.Range("B42").Value = 'spinnen'
$ws.Range("B43").Value = 'quälen'
$ws.Range("B44").Value = 'folgen'
$ws.Range("B45").Value = 'füllen'
$ws.Range("B46").Value = 'liefern'
$ws.Range("B47").Value = 'graben'
$ws.Range("B48").Value = 'planen'
$ws.Range("B49").Value = 'binden'
$ws.Range("B50").Value = 'sparen'
$ws.Range("B51").Value = 'gelten'
$ws.Range("B52").Value = 'zünden'
$ws.Range("B53").Value = 'bluten'
$ws.Range("B54").Value = 'wählen'
$ws.Range("B55").Value = 'werden'
$ws.Range("B56").Value = 'wüten'
$ws.Range("B57").Value = 'lockern'
$ws.Range("B58").Value = 'räumen'
$ws.Range("B59").Value = 'zählen'
$ws.Range("B60").Value = 'mauern'
$ws.Range("B61").Value = 'seufzen'
$ws.Range("B62").Value = 'steuern'
$ws.Range("B63").Value = 'knurren'
$ws.Range("B64").Value = 'erben'
$ws.Range("B65").Value = 'deuten'
$ws.Range("B66").Value = 'schrecken'
$ws.Range("B67").Value = 'beten'
$ws.Range("B68").Value = 'öffnen'
$ws.Range("B69").Value = 'fallen'
$ws.Range("B70").Value = 'malen'
$ws.Range("B71").Value = 'starten'
$ws.Range("B72").Value = 'treffen'
$ws.Range("B73").Value = 'leiden'
$ws.Range("B74").Value = 'reizen'
$ws.Range("B75").Value = 'äußern'
$ws.Range("B76").Value = 'bellen'
$ws.Range("B77").Value = 'stehlen'
$ws.Range("B78").Value = 'fügen'
$ws.Range("B79").Value = 'knarren'
$ws.Range("B80").Value = 'leeren'
$ws.Range("B81").Value = 'rühren'
$ws.Range("B82").Value = 'führen'
$ws.Range("B83").Value = 'ächzen'
$ws.Range("B84").Value = 'kaufen'
$ws.Range("B85").Value = 'teilen'
$ws.Range("B86").Value = 'fahren'
$ws.Range("B87").Value = 'eignen'
$ws.Range("B88").Value = 'testen'
$ws.Range("B89").Value = 'sagen'
$ws.Range("B90").Value = 'altern'
$ws.Range("B91").Value = 'spannen'
$ws.Range("B92").Value = 'wachsen'
$ws.Range("B93").Value = 'gleichen'
$ws.Range("B94").Value = 'mögen'
$ws.Range("B95").Value = 'pfeifen'
$ws.Range("B96").Value = 'kratzen'
$ws.Range("B97").Value = 'rufen'
$ws.Range("B98").Value = 'beißen'
$ws.Range("B99").Value = 'trauen'
$ws.Range("B100").Value = 'tropfen'
$ws.Range("B101").Value = 'regeln'
$ws.Range("B102").Value = 'kümmern'
$ws.Range("B103").Value = 'dürfen'
$ws.Range("B104").Value = 'warnen'
$ws.Range("B105").Value = 'stoßen'
$ws.Range("B106").Value = 'wirken'
$ws.Range("B107").Value = 'heulen'
$ws.Range("B108").Value = 'flüstern'
$ws.Range("B109").Value = 'heilen'
$ws.Range("B110").Value = 'bremsen'
$ws.Range("B111").Value = 'nähen'
$ws.Range("B112").Value = 'treiben'
$ws.Range("B113").Value = 'suchen'
$ws.Range("B114").Value = 'sichern'
$ws.Range("B115").Value = 'schreiten'
$ws.Range("B116").Value = 'stecken'
$ws.Range("B117").Value = 'schütteln'
$ws.Range("B118").Value = 'bauen'
$ws.Range("B119").Value = 'heben'
$ws.Range("B120").Value = 'saufen'
$ws.Range("B121").Value = 'pflanzen'
$ws.Range("B122").Value = 'wehtun'
$ws.Range("B123").Value = 'kichern'
$ws.Range("B124").Value = 'schützen'
$ws.Range("B125").Value = 'stammen'
$ws.Range("B126").Value = 'irren'
$ws.Range("B127").Value = 'brauchen'
$ws.Range("B128").Value = 'boxen'
$ws.Range("B129").Value = 'stürmen'
$ws.Range("B130").Value = 'fangen'
$ws.Range("B131").Value = 'wundern'
$ws.Range("B132").Value = 'feiern'
$ws.Range("B133").Value = 'handeln'
$ws.Range("B134").Value = 'starren'
$ws.Range("B135").Value = 'fällen'
$ws.Range("B136").Value = 'klingen'
$ws.Range("B137").Value = 'decken'
$ws.Range("B138").Value = 'sammeln'
$ws.Range("B139").Value = 'kosten'
$ws.Range("B140").Value = 'tollen'
$ws.Range("B141").Value = 'süßen'
$ws.Range("B142").Value = 'schwingen'
$ws.Range("B143").Value = 'töten'
$ws.Range("B144").Value = 'rasen'
$ws.Range("B145").Value = 'proben'
$ws.Range("B146").Value = 'schlucken'
$ws.Range("B147").Value = 'stören'
$ws.Range("B148").Value = 'ehren'
$ws.Range("B149").Value = 'zeigen'
$ws.Range("B150").Value = 'sinken'
$ws.Range("B151").Value = 'kranken'
$ws.Range("B152").Value = 'formen'
$ws.Range("B153").Value = 'bergen'
$ws.Range("B154").Value = 'wecken'
$ws.Range("B155").Value = 'streichen'
$ws.Range("B156").Value = 'achten'
$ws.Range("B157").Value = 'sterben'
$ws.Range("B158").Value = 'spenden'
$ws.Range("B159").Value = 'schmecken'
$ws.Range("B160").Value = 'zielen'
$ws.Range("B161").Value = 'grüßen'
$ws.Range("B162").Value = 'füttern'
$ws.Range("B163").Value = 'wenden'
$ws.Range("B164").Value = 'baden'
$ws.Range("B165").Value = 'dringen'
$ws.Range("B166").Value = 'heißen'
$ws.Range("B167").Value = 'liegen'
$ws.Range("B168").Value = 'sorgen'
$ws.Range("B169").Value = 'orten'
$ws.Range("B170").Value = 'trotzen'
$ws.Range("B171").Value = 'ärgern'
$ws.Range("B172").Value = 'fischen'
$ws.Range("B173").Value = 'fordern'
$ws.Range("B174").Value = 'gründen'
$ws.Range("B175").Value = 'wandern'
$ws.Range("B176").Value = 'lügen'
$ws.Range("B177").Value = 'helfen'
$ws.Range("B178").Value = 'hauen'
$ws.Range("B179").Value = 'jagen'
$ws.Range("B180").Value = 'scheinen'
$ws.Range("B181").Value = 'stillen'
$ws.Range("B182").Value = 'arten'
$ws.Range("B183").Value = 'platzen'
$ws.Range("B184").Value = 'dienen'
$ws.Range("B185").Value = 'filmen'
$ws.Range("B186").Value = 'werfen'
$ws.Range("B187").Value = 'parken'
$ws.Range("B188").Value = 'drehen'
$ws.Range("B189").Value = 'backen'
$ws.Range("B190").Value = 'schenken'
$ws.Range("B191").Value = 'klettern'
$ws.Range("B192").Value = 'buchen'
$ws.Range("B193").Value = 'ändern'
